# Insert a new row of data at row 2 (pushing the existing rows 2 and 3
# down to rows 3 and 4), then populate the new row 2 with the new
# invigilator's details.
#
# We deliberately avoid Range.Value as a *read* source (this runtime's
# COM shim doesn't surface it reliably) and instead use Range.Copy(dest),
# which duplicates both value and the exact existing style index without
# fabricating new style entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows down: row 3 -> row 4, then row 2 -> row 3.
# (Do the lower row first so we don't clobber it.)
$ws.Range("A3:F3").Copy($ws.Range("A4:F4"))
$ws.Range("A2:F2").Copy($ws.Range("A3:F3"))

# Fill in the new row 2 with the newly added invigilator's info.
$ws.Range("A2").Value = "wrc teacher``"
$ws.Range("B2").Value = "test@gmail.com"

# Phone number must stay literal text (keep the leading "+"); a leading
# apostrophe forces text entry, then we reset the style back to Normal so
# we don't leave a stray quote-prefixed style applied to the cell.
$ws.Range("C2").Value = "'+13442322"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "Gagangauda"
$ws.Range("E2").Value = "Department of Computer and Electronics Engineering"
$ws.Range("F2").Value = "Professor"
